$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: clear out the Species (perc.) column (C) and
#     zero out the Species (no.) column (B) for rows 2-7 -------------------
$ws = $wb.Worksheets.Item("Range Status")
$ws.Range("B2:B7").Value = 0
$ws.Range("C2:C7").ClearContents()

# --- "Species qualification" sheet: Range Analysis row (row 5) now has 0
#     species selected ------------------------------------------------------
$ws = $wb.Worksheets.Item("Species qualification")
$ws.Range("B5").Value = 0

# --- "High Priority break-up" sheet: add the "New High Species" counts for
#     the IUCN row (row 2) ---------------------------------------------------
$ws = $wb.Worksheets.Item("High Priority break-up")
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 100
